$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2:A11 values from 4 to 1
$ws.Range("A2:A11").Value = 1

# Update the selection to match the diff (active cell A2, selection A2:A11)
$ws.Range("A2:A11").Select()
